$wb = $excel.ActiveWorkbook

# Update population figures for France (rows 4 and 5) on the "pop" sheet
# to include revised 2014/2015 values (as part of refreshed demography data).
$popSheet = $wb.Worksheets.Item("pop")
$popSheet.Range("D4").Value = 32045129
$popSheet.Range("E4").Value = 32174258
$popSheet.Range("D5").Value = 34120851
$popSheet.Range("E5").Value = 34283895

# Remove the obsolete "__groups__" sheet (even_years@time / odd_years@time)
$groupsSheet = $wb.Worksheets.Item("__groups__")
$groupsSheet.Delete()

# Make "__axes__" the active/selected sheet
$axesSheet = $wb.Worksheets.Item("__axes__")
$axesSheet.Activate()
